# Generate Report for Archive
# The localization status for this item moved from "Ready for handoff" to
# "In Translation". That status string shows up in the per-language Status
# column on the "zh-cn" / "de-de" sheets, and is mirrored on the "Overview"
# sheet (one column per language).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn!C2 and de-de!C2 ("Status" column)
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# The Status column was sized to fit its text ("Ready for handoff"); now that
# the text is shorter ("In Translation") re-fit those columns so they shrink
# to match the new (shorter) status text.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
